$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet
$ws.Name = "ChucnySlot RUF Gen algs"

# Add credit/info text below the existing data
$ws.Range("B24").Value = "These are all the 48 ChucnySlot EOCPLS insertion algs. "
$ws.Range("B25").Value = "Made by: Chucny"
$ws.Range("B26").Value = "DM @chucny on Discord if you have any questions"
$ws.Range("B27").Value = "Search ""ChucnyZB"" on speedsolving wiki if you want full documentation of the method."

$ws.Range("B27").Select() | Out-Null
